$d = $word.ActiveDocument

# Helper: replace $oldText with $newText while preserving a leading empty
# run (a bare <w:r/>) that sits immediately before the text-bearing run in
# the original document. Plain Find/Replace (or a Range.Text assignment)
# causes this engine to defragment/merge the paragraph's runs and silently
# drop the empty run, so instead we (1) insert the new text as a brand new
# run right before the old text [this is the one operation observed not to
# trigger the defrag pass], (2) re-apply bold/italic formatting to the new
# run if the old run had it, and (3) delete the old text afterwards.
function Replace-KeepRun {
    param(
        [string]$OldText,
        [string]$NewText
    )

    $rng = $d.Content
    $found = $rng.Find.Execute($OldText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Output "NOT FOUND: $OldText"
        return
    }

    $isBold = $rng.Bold
    $isItalic = $rng.Italic
    $startPos = $rng.Start
    $endPos = $rng.End

    $rng.InsertBefore($NewText)

    $newRng = $d.Range($startPos, $startPos + $NewText.Length)
    if ($isBold) { $newRng.Bold = 1 }
    if ($isItalic) { $newRng.Italic = 1 }

    $shift = $NewText.Length
    $delRng = $d.Range($startPos + $shift, $endPos + $shift)
    $delRng.Delete()
}

# Title (Heading1) -- no leading empty run here, plain replace is safe.
$d.Content.Find.Execute("Play Electric Avenue for Free | 80s-Themed Slot Game", $true, $false, $false, $false, $false, $true, 1, $false, "Play Electric Avenue Free: Review of Exciting Slot Game", 2)

# "What we like" bullet list (each has a leading empty run to preserve)
Replace-KeepRun "Two modes of free spins with various benefits" "Exciting gameplay features with Wild Reels and multipliers"
Replace-KeepRun "High winning potential with a maximum payout of 6,833x" "Two modes of free spins with different bonus features"
Replace-KeepRun "Appealing 80s theme with neon symbols and electrifying music" "Flamboyant 80s theme with bright neon symbols"
Replace-KeepRun "Seamless mobile compatibility for on-the-go play" "High volatility for players looking for high payout potential"

# "What we don't like" bullet list
Replace-KeepRun "High volatility may not be suitable for low-budget players" "Maximum payout potential not yet released"
Replace-KeepRun "The maximum payout potential is yet to be released. " "Limited bet range from 20 cents to 30 euros per spin"

# Bold repeated title near the bottom
Replace-KeepRun "Play Electric Avenue for Free | 80s-Themed Slot Game" "Play Electric Avenue Free: Review of Exciting Slot Game"

# Italic meta description
Replace-KeepRun "Read our review of Electric Avenue, a 6-reel, 4096 ways to win slot with Wild Reels and multipliers. Play for free and experience the neon 80s-inspired slot game." "Read our review of Electric Avenue and play for free to experience its thrilling gameplay features"
